$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("A2").Value = 7474
$ws.Range("B2").Value = "Agatha Sales"
$ws.Range("C2").Value = "Marketing"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 45106
$ws.Range("G2").Value = 11163.84

# Row 3
$ws.Range("A3").Value = 87852
$ws.Range("B3").Value = "Sophia Almeida"
$ws.Range("C3").Value = "Jurídico"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 8
$ws.Range("G3").Value = 3294.68

# Row 4
$ws.Range("A4").Value = 84057
$ws.Range("B4").Value = "Vicente Vieira"
$ws.Range("C4").Value = "Financeiro"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 45101
$ws.Range("G4").Value = 11375.62

# Row 5
$ws.Range("A5").Value = 42868
$ws.Range("B5").Value = "Theo Caldeira"
$ws.Range("C5").Value = "TI"
$ws.Range("E5").Value = 6
$ws.Range("G5").Value = 7216.93

# Row 6
$ws.Range("A6").Value = 73610
$ws.Range("B6").Value = "Camila Pereira"
$ws.Range("C6").Value = "P&D"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45101
$ws.Range("G6").Value = 8403.92

# Row 7
$ws.Range("A7").Value = 60732
$ws.Range("B7").Value = "Yasmin Pires"
$ws.Range("C7").Value = "Financeiro"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 7
$ws.Range("F7").Value = 45083
$ws.Range("G7").Value = 2999.94

# Row 8
$ws.Range("A8").Value = 10123
$ws.Range("B8").Value = "Luiza Cavalcanti"
$ws.Range("C8").Value = "Vendas"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 7
$ws.Range("F8").Value = 45094
$ws.Range("G8").Value = 6239.67

# Row 9
$ws.Range("A9").Value = 87084
$ws.Range("B9").Value = "João Guilherme Cunha"
$ws.Range("C9").Value = "Vendas"
$ws.Range("D9").Value = "Doença"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 45084
$ws.Range("G9").Value = 12099.85

# Row 10
$ws.Range("A10").Value = 30268
$ws.Range("B10").Value = "Thales da Rocha"
$ws.Range("C10").Value = "Operações"
$ws.Range("D10").Value = "Outros"
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 45099
$ws.Range("G10").Value = 10832.79

# Row 11
$ws.Range("A11").Value = 22946
$ws.Range("B11").Value = "Sr. Emanuel Viana"
$ws.Range("C11").Value = "Jurídico"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 5
$ws.Range("G11").Value = 5418.21
